$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellRuns($cell, $xmlRunsInner) {
    $rng = $cell.Range
    # Exclude the trailing cell-end mark so InsertXML replaces the
    # paragraph's content in place instead of appending a sibling paragraph.
    $target = $d.Range($rng.Start, $rng.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $xmlRunsInner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# 1) "+1 pt" -> "+2 pt" (same visual result, split into 3 runs: "+", "2", " pt")
$cellBonus = $t.Rows.Item(7).Cells.Item(1)
Set-CellRuns $cellBonus '<w:r><w:t>+</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> pt</w:t></w:r>'

# 2) The two standalone "+0.5 pt" cells (rows 9 and 10, first column) join the
#    "+2 pt" vertical-merge group started at row 7 -> become vMerge continuation
#    cells with empty paragraphs, same as row 8 already is.
$t.Rows.Item(7).Cells.Item(1).Merge($t.Rows.Item(10).Cells.Item(1))

# 3) "Total (11 pts possibles)" -> "Total (12 pts possibles)" (split into 3 runs)
$cellTotal = $t.Rows.Item(11).Cells.Item(2)
Set-CellRuns $cellTotal '<w:r><w:t>Total (1</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> pts possibles)</w:t></w:r>'
